$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.527.32'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.922.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.42%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4819'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4067'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("E10").Value = '  -0.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.905.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.077'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.271'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.74'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06872'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.012'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001038'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.561.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.684'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.183'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.126.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.02'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.485'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.099'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.016'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09639'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.627'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("E34").Value = '  -0.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.376'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06370'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.192'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5960'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.883'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1848'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.465'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.279'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07502'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5572'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("E48").Value = '  +0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.437'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '
